$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2025-08-31"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = 56.43000030517578
$ws.Range("C16").Value = 669
$ws.Range("D16").Value = 313.9500122070312
